# Updates R programs and scenarios
# Fill in row 11 (2017 Q4) values for columns C:I (GDP, UEMP, CPI, LTRate, EURUSD, WTI, RPP)
# on both "Test 1" and "Test 2" worksheets.

$wb = $excel.ActiveWorkbook

$values = @{
    "C11" = 0.6462729991457934
    "D11" = -0.20000000000000018
    "E11" = 1.222806544026625
    "F11" = 0.03760000000000008
    "G11" = 1.5829618029997903
    "H11" = 16.12947350163202
    "I11" = 1.29289423366793
}

foreach ($sheetName in @("Test 1", "Test 2")) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($addr in $values.Keys) {
        $ws.Range($addr).Value = $values[$addr]
    }
}
